# Updates cryptos list data (prices / volume %) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text entry (matches the source inlineStr cells)
# so values like '229.72' or '1.003' are not auto-coerced to numbers by Excel,
# and the cell's NumberFormat stays General (unchanged) rather than becoming '@'.
$ws.Range('D2').Value = "'28.005.73"
$ws.Range('E2').Value = "'  +6.64%  "
$ws.Range('D3').Value = "'1.741.01"
$ws.Range('E3').Value = "'  +5.08%  "
$ws.Range('E4').Value = "'  -0.15%  "
$ws.Range('D5').Value = "'229.72"
$ws.Range('E5').Value = "'  +4.71%  "
$ws.Range('D6').Value = "'0.5441"
$ws.Range('E6').Value = "'  +3.90%  "
$ws.Range('E7').Value = "'  -0.19%  "
$ws.Range('D8').Value = "'0.2786"
$ws.Range('E8').Value = "'  +4.30%  "
$ws.Range('D9').Value = "'0.06725"
$ws.Range('E9').Value = "'  +5.59%  "
$ws.Range('D10').Value = "'21.86"
$ws.Range('E10').Value = "'  +5.60%  "
$ws.Range('D11').Value = "'0.07802"
$ws.Range('E11').Value = "'  +1.32%  "
$ws.Range('D12').Value = "'4.714"
$ws.Range('D13').Value = "'1.744.34"
$ws.Range('E13').Value = "'  +11.03%  "
$ws.Range('D14').Value = "'1.976.35"
$ws.Range('E14').Value = "'  +4.83%  "
$ws.Range('D15').Value = "'0.6020"
$ws.Range('E15').Value = "'  +6.54%  "
$ws.Range('D16').Value = "'0.0₅8436"
$ws.Range('E16').Value = "'  +1.86%  "
$ws.Range('D17').Value = "'69.69"
$ws.Range('E17').Value = "'  +6.43%  "
$ws.Range('D18').Value = "'27.978.44"
$ws.Range('E18').Value = "'  +6.54%  "
$ws.Range('D19').Value = "'226.37"
$ws.Range('E19').Value = "'  +17.31%  "
$ws.Range('D20').Value = "'4.836"
$ws.Range('E20').Value = "'  +2.97%  "
$ws.Range('E21').Value = "'  -0.18%  "
$ws.Range('D22').Value = "'10.97"
$ws.Range('E22').Value = "'  +4.92%  "
$ws.Range('D23').Value = "'6.281"
$ws.Range('E23').Value = "'  +4.51%  "
$ws.Range('E24').Value = "'  -0.23%  "
$ws.Range('D25').Value = "'147.33"
$ws.Range('E25').Value = "'  +2.78%  "
$ws.Range('D26').Value = "'0.1252"
$ws.Range('E26').Value = "'  +4.20%  "
$ws.Range('D27').Value = "'7.464"
$ws.Range('E27').Value = "'  +2.33%  "
$ws.Range('D28').Value = "'17.03"
$ws.Range('E28').Value = "'  +6.91%  "
$ws.Range('D29').Value = "'1.636"
$ws.Range('E29').Value = "'  +8.96%  "
$ws.Range('D30').Value = "'0.05635"
$ws.Range('E30').Value = "'  -0.15%  "
$ws.Range('E31').Value = "'  +3.49%  "
$ws.Range('D32').Value = "'3.721"
$ws.Range('E32').Value = "'  +5.97%  "
$ws.Range('D33').Value = "'3.552"
$ws.Range('E33').Value = "'  +5.86%  "
$ws.Range('D34').Value = "'1.659"
$ws.Range('E34').Value = "'  +4.70%  "
$ws.Range('D35').Value = "'0.9866"
$ws.Range('E35').Value = "'  +4.07%  "
$ws.Range('D36').Value = "'2.857"
$ws.Range('E36').Value = "'  +1.82%  "
$ws.Range('D37').Value = "'2.449"
$ws.Range('E37').Value = "'  +1.41%  "
$ws.Range('D38').Value = "'0.5947"
$ws.Range('E38').Value = "'  +3.06%  "
$ws.Range('D39').Value = "'0.01684"
$ws.Range('E39').Value = "'  +5.11%  "
$ws.Range('D40').Value = "'5.959"
$ws.Range('E40').Value = "'  -0.27%  "
$ws.Range('B41').Value = "'TrustWalletToken"
$ws.Range('C41').Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('D41').Value = "'0.8477"
$ws.Range('E41').Value = "'  +0.20%  "
$ws.Range('B42').Value = "'Maker"
$ws.Range('C42').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D42').Value = "'1.046.76"
$ws.Range('E42').Value = "'  +2.88%  "
$ws.Range('E43').Value = "'  -0.14%  "
$ws.Range('D44').Value = "'102.47"
$ws.Range('E44').Value = "'  +0.62%  "
$ws.Range('D45').Value = "'1.880.83"
$ws.Range('E45').Value = "'  +4.69%  "
$ws.Range('E46').Value = "'  +10.75%  "
$ws.Range('D47').Value = "'60.29"
$ws.Range('E47').Value = "'  +3.23%  "
$ws.Range('D48').Value = "'8.301"
$ws.Range('E48').Value = "'  +3.39%  "
$ws.Range('D49').Value = "'0.4421"
$ws.Range('E49').Value = "'  +1.67%  "
$ws.Range('D50').Value = "'1.003"
$ws.Range('E50').Value = "'  +0.05%  "
$ws.Range('D51').Value = "'0.05309"
$ws.Range('E51').Value = "'  -0.26%  "
